# Applies a weekly refresh of the Jengibre (Arica) price data: the rows of the
# data table (rows 2-32) get re-shuffled to a new order, changing the values
# in columns D (Fecha), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), N (Unidad de
# comercializacion), P (Precio $/Kg) and Q (Kg o Unidades) while A, B, C, E,
# F, G, H, O and R stay constant for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> old row number that supplies its new content.
$rowMap = @{
    2  = 28
    3  = 2
    4  = 7
    5  = 25
    6  = 17
    7  = 11
    8  = 27
    9  = 18
    10 = 16
    11 = 31
    12 = 22
    13 = 5
    14 = 3
    15 = 15
    16 = 6
    17 = 29
    18 = 14
    19 = 13
    20 = 19
    21 = 32
    22 = 9
    23 = 20
    24 = 24
    25 = 4
    26 = 12
    27 = 10
    28 = 23
    29 = 21
    30 = 26
    31 = 30
    32 = 8
}

$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

# Snapshot the original values of every relevant column for every row before
# any writes happen, since the remap reorders rows (reading after writing
# would use already-overwritten data).
$snapshot = @{}
for ($r = 2; $r -le 32; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the new values into each row according to the mapping.
foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $srcVals = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $srcVals[$c]
    }
}
